# Updated cryptos list on Sun Feb 26 23:30:00 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.501.17"
$ws.Range("E2").Value = "  +1.72%  "
$ws.Range("D3").Value = "1.639.98"
$ws.Range("E3").Value = "  +3.06%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").Value = "'308.44"
$ws.Range("E5").Value = "  +2.28%  "
$ws.Range("D6").Value = "'0.9996"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").Value = "'0.3770"
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").Value = "'52.86"
$ws.Range("E8").Value = "  +3.45%  "
$ws.Range("D9").Value = "'0.3683"
$ws.Range("E9").Value = "  +2.05%  "
$ws.Range("D10").Value = "'1.277"
$ws.Range("E10").Value = "  +2.54%  "
$ws.Range("D11").Value = "'0.08211"
$ws.Range("E11").Value = "  +2.13%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("D13").Value = "'23.21"
$ws.Range("E13").Value = "  +3.54%  "
$ws.Range("D14").Value = "'6.668"
$ws.Range("E14").Value = "  +2.27%  "
$ws.Range("D15").Value = "'0.00001281"
$ws.Range("E15").Value = "  +3.52%  "
$ws.Range("D16").Value = "'7.469"
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("D17").Value = "1.639.94"
$ws.Range("E17").Value = "  +2.67%  "
$ws.Range("D18").Value = "'95.08"
$ws.Range("E18").Value = "  +2.32%  "
$ws.Range("D19").Value = "'0.06961"
$ws.Range("E19").Value = "  +2.90%  "
$ws.Range("E20").Value = "  +2.73%  "
$ws.Range("D21").Value = "'6.585"
$ws.Range("E21").Value = "  +2.18%  "
$ws.Range("D22").Value = "'0.9968"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "23.495.46"
$ws.Range("E23").Value = "  +1.74%  "
$ws.Range("D24").Value = "'12.96"
$ws.Range("E24").Value = "  +1.19%  "
$ws.Range("D25").Value = "'3.116"
$ws.Range("E25").Value = "  +6.59%  "
$ws.Range("D26").Value = "'2.410"
$ws.Range("E26").Value = "  +1.14%  "
$ws.Range("D27").Value = "'21.43"
$ws.Range("E27").Value = "  +2.51%  "
$ws.Range("D28").Value = "'151.56"
$ws.Range("E28").Value = "  +2.06%  "
$ws.Range("D29").Value = "'5.334"
$ws.Range("E29").Value = "  +2.85%  "
$ws.Range("D30").Value = "'136.14"
$ws.Range("E30").Value = "  +2.38%  "
$ws.Range("D31").Value = "'2.421"
$ws.Range("E31").Value = "  +2.07%  "
$ws.Range("D32").Value = "'6.866"
$ws.Range("E32").Value = "  +2.31%  "
$ws.Range("D33").Value = "1.820.58"
$ws.Range("E33").Value = "  +3.06%  "
$ws.Range("D34").Value = "'0.9776"
$ws.Range("E34").Value = "  +2.13%  "
$ws.Range("D35").Value = "'0.02815"
$ws.Range("E35").Value = "  +5.33%  "
$ws.Range("D36").Value = "'10.46"
$ws.Range("E36").Value = "  +3.58%  "
$ws.Range("D37").Value = "'0.07476"
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").Value = "'6.222"
$ws.Range("E38").Value = "  +2.11%  "
$ws.Range("D39").Value = "'0.2542"
$ws.Range("E39").Value = "  +1.74%  "
$ws.Range("D40").Value = "'0.08879"
$ws.Range("E40").Value = "  +1.13%  "
$ws.Range("D41").Value = "'1.403"
$ws.Range("E41").Value = "  +3.41%  "
$ws.Range("D42").Value = "'0.7164"
$ws.Range("E42").Value = "  +0.99%  "
$ws.Range("D43").Value = "'12.65"
$ws.Range("E43").Value = "  +3.86%  "
$ws.Range("D44").Value = "'16.17"
$ws.Range("E44").Value = "  +8.16%  "
$ws.Range("D45").Value = "'0.6621"
$ws.Range("E45").Value = "  +2.03%  "
$ws.Range("D46").Value = "'2.362"
$ws.Range("E46").Value = "  +3.62%  "
$ws.Range("D47").Value = "'4.046"
$ws.Range("E47").Value = "  +1.37%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.08068"
$ws.Range("E48").Value = "  +2.08%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'130.90"
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("B50").Value = "Flow"
$ws.Range("C50").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D50").Value = "'1.221"
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "'1.218"
$ws.Range("E51").Value = "  +2.68%  "